$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A28").Value = 28
$ws.Range("B28").Value = "ds33"
$ws.Range("C28").Value = "default.jpg"
$ws.Range("D28").Value = "Ngô Xuân Hinh"
$ws.Range("E28").Value = "Nam"
# F28 ("Ngày sinh") must stay as the literal text "2022-09-16" instead of being
# auto-converted into a date serial number. Copy an existing cell that already
# holds this exact text (as a shared string, style 0) rather than assigning
# the string via .Value, which would trigger Excel's date auto-detection.
$ws.Range("F8").Copy($ws.Range("F28"))
$ws.Range("G28").Value = "1,2"
$ws.Range("H28").Value = "N/A"
$ws.Range("I28").Value = "N/A"
$ws.Range("J28").Value = "N/A"
$ws.Range("K28").Value = "N/A"
$ws.Range("L28").Value = "N/A"
$ws.Range("M28").Value = "Chưa có"
$ws.Range("N28").Value = "Chưa có"
$ws.Range("O28").Value = $false
$ws.Range("P28").Value = $false
